$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing rows (2..13) down to (3..14)
$ws.Rows("2:2").Insert()

# Populate the new row 2 with the NFE47-TP-REG field definition
$ws.Cells.Item(2, 1).Value = "NFE47-TP-REG"
$ws.Cells.Item(2, 2).Value = 1
$ws.Cells.Item(2, 3).Value = 2
$ws.Cells.Item(2, 4).Value = "NUMERO"
$ws.Cells.Item(2, 5).Value = "N"
$ws.Cells.Item(2, 6).Value = "'"

# New rows inherit formatting from the row above; strip that back to the
# unstyled look used by the other data rows.
$ws.Range("A2:F2").ClearFormats()

$wb.Save()
